$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 38325.11
$ws.Range("I38").Value = 41865.875
$ws.Range("J38").Value = 9999
$ws.Range("K38").Value = 125597.625
$ws.Range("L38").Value = 29997
$ws.Range("M38").Value = -125225.625
$ws.Range("N38").Value = -30741

$ws.Range("H43").Value = 6907.636
$ws.Range("J43").Value = 6997.5
$ws.Range("L43").Value = 6997.5
$ws.Range("N43").Value = -7135.5

$ws.Range("H98").Value = 15625908
$ws.Range("I98").Value = 17857780
$ws.Range("K98").Value = 17857780
$ws.Range("M98").Value = -17856282

$ws.Range("H103").Value = 38462980
$ws.Range("J103").Value = 62501816
$ws.Range("L103").Value = 187505448
$ws.Range("N103").Value = -187506620

$ws.Range("H122").Value = 15625908
$ws.Range("I122").Value = 17857780
$ws.Range("K122").Value = 53573340
$ws.Range("M122").Value = -53570890

$ws.Range("H135").Value = 1590
$ws.Range("I135").Value = 288.33334
$ws.Range("K135").Value = 2595.00006
$ws.Range("M135").Value = -60.0000600000003

$ws.Range("H141").Value = 5927.2144
$ws.Range("I141").Value = 5844.846
$ws.Range("K141").Value = 17534.538
$ws.Range("M141").Value = -12354.538

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 71431704
$ws.Range("I4").Value = 3920.6
$ws.Range("J4").Value = 250001170
$ws.Range("K4").Value = 3920.6
$ws.Range("L4").Value = 250001170
$ws.Range("M4").Value = -3804.6
$ws.Range("N4").Value = -250001402

$ws.Range("H58").Value = 750022500
$ws.Range("J58").Value = 750022500
$ws.Range("L58").Value = 750022500
$ws.Range("N58").Value = -750023360

$ws.Range("H63").Value = 3899.4
$ws.Range("I63").Value = 4099
$ws.Range("K63").Value = 4099
$ws.Range("M63").Value = -3413

$ws.Range("H66").Value = 3899.4
$ws.Range("I66").Value = 4099
$ws.Range("K66").Value = 20495
$ws.Range("M66").Value = -17063

$ws.Range("H97").Value = 2595.2856
$ws.Range("I97").Value = 1861
$ws.Range("J97").Value = 3329.5715
$ws.Range("K97").Value = 1861
$ws.Range("L97").Value = 3329.5715
$ws.Range("M97").Value = -1365
$ws.Range("N97").Value = -4321.5715

$ws.Range("H102").Value = 3413
$ws.Range("I102").Value = 3413
$ws.Range("K102").Value = 3413
$ws.Range("M102").Value = -1791

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1948.3928
$ws.Range("I94").Value = 2303.1052
$ws.Range("J94").Value = 1199.5555
$ws.Range("K94").Value = 2303.1052
$ws.Range("L94").Value = 1199.5555
$ws.Range("M94").Value = -1852.1052
$ws.Range("N94").Value = -2101.5555

$ws.Range("H134").Value = 3228213.5
$ws.Range("I134").Value = 2480.28
$ws.Range("J134").Value = 16668769
$ws.Range("K134").Value = 7440.84
$ws.Range("L134").Value = 50006307
$ws.Range("M134").Value = -4905.84
$ws.Range("N134").Value = -50011377

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3827.0715
$ws.Range("I58").Value = 2756.7
$ws.Range("J58").Value = 6503
$ws.Range("K58").Value = 2756.7
$ws.Range("L58").Value = 6503
$ws.Range("M58").Value = -2553.7
$ws.Range("N58").Value = -6909

$ws.Range("H99").Value = 14200.7
$ws.Range("I99").Value = 8065
$ws.Range("K99").Value = 8065
$ws.Range("M99").Value = -6567

$ws.Range("H107").Value = 1291.7826
$ws.Range("I107").Value = 419.1875
$ws.Range("J107").Value = 3286.2856
$ws.Range("K107").Value = 419.1875
$ws.Range("L107").Value = 3286.2856
$ws.Range("M107").Value = 1500.8125
$ws.Range("N107").Value = -7126.2856

$ws.Range("H111").Value = 75175
$ws.Range("J111").Value = 75175
$ws.Range("L111").Value = 75175
$ws.Range("N111").Value = -83355

$ws.Range("H122").Value = 1755.409
$ws.Range("J122").Value = 2795
$ws.Range("L122").Value = 8385
$ws.Range("N122").Value = -13285

$ws.Range("H126").Value = 14200.7
$ws.Range("I126").Value = 8065
$ws.Range("K126").Value = 24195
$ws.Range("M126").Value = -21725

$ws.Range("H132").Value = 3479.1177
$ws.Range("I132").Value = 2405
$ws.Range("J132").Value = 4687.5
$ws.Range("K132").Value = 7215
$ws.Range("L132").Value = 14062.5
$ws.Range("M132").Value = -4685
$ws.Range("N132").Value = -19122.5

$ws.Range("H134").Value = 2137.3845
$ws.Range("I134").Value = 1707.7368
$ws.Range("K134").Value = 5123.2104
$ws.Range("M134").Value = -2588.2104

$ws.Range("H136").Value = 3827.0715
$ws.Range("I136").Value = 2756.7
$ws.Range("J136").Value = 6503
$ws.Range("K136").Value = 8270.099999999999
$ws.Range("L136").Value = 19509
$ws.Range("M136").Value = -5720.099999999999
$ws.Range("N136").Value = -24609

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 2495.1667
$ws.Range("J23").Value = 2495.1667
$ws.Range("L23").Value = 7485.500100000001
$ws.Range("N23").Value = -7955.500100000001

$ws.Range("H50").Value = 159.57143
$ws.Range("I50").Value = 153.6
$ws.Range("J50").Value = 174.5
$ws.Range("K50").Value = 460.8
$ws.Range("L50").Value = 523.5
$ws.Range("M50").Value = 20.20000000000005
$ws.Range("N50").Value = -1485.5

$ws.Range("H53").Value = 159.57143
$ws.Range("I53").Value = 153.6
$ws.Range("J53").Value = 174.5
$ws.Range("K53").Value = 460.8
$ws.Range("L53").Value = 523.5
$ws.Range("M53").Value = 20.20000000000005
$ws.Range("N53").Value = -1485.5

$ws.Range("H117").Value = 7688.4287
$ws.Range("I117").Value = 712.5
$ws.Range("J117").Value = 10478.8
$ws.Range("K117").Value = 2137.5
$ws.Range("L117").Value = 31436.4
$ws.Range("M117").Value = 1304.5
$ws.Range("N117").Value = -38320.39999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 35318.816
$ws.Range("I58").Value = 27056.445
$ws.Range("J58").Value = 72499.5
$ws.Range("K58").Value = 27056.445
$ws.Range("L58").Value = 72499.5
$ws.Range("M58").Value = -26779.445
$ws.Range("N58").Value = -73053.5

$ws.Range("H102").Value = 1782.3871
$ws.Range("I102").Value = 1814.4286
$ws.Range("K102").Value = 1814.4286
$ws.Range("M102").Value = -192.4286

$ws.Range("H105").Value = 59417.25
$ws.Range("I105").Value = 27000
$ws.Range("J105").Value = 70223
$ws.Range("K105").Value = 27000
$ws.Range("L105").Value = 70223
$ws.Range("M105").Value = -23506
$ws.Range("N105").Value = -77211

$ws.Range("H132").Value = 2176849.2
$ws.Range("I132").Value = 2961.7
$ws.Range("J132").Value = 16669433
$ws.Range("K132").Value = 8885.099999999999
$ws.Range("L132").Value = 50008299
$ws.Range("M132").Value = -6355.099999999999
$ws.Range("N132").Value = -50013359

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2460
$ws.Range("I16").Value = 2631.6155
$ws.Range("K16").Value = 2631.6155
$ws.Range("M16").Value = -2461.6155

$ws.Range("H22").Value = 2855.3333
$ws.Range("I22").Value = 1100
$ws.Range("J22").Value = 4259.6
$ws.Range("K22").Value = 1100
$ws.Range("L22").Value = 4259.6
$ws.Range("M22").Value = -805
$ws.Range("N22").Value = -4849.6

$ws.Range("H27").Value = 2855.3333
$ws.Range("I27").Value = 1100
$ws.Range("J27").Value = 4259.6
$ws.Range("K27").Value = 1100
$ws.Range("L27").Value = 4259.6
$ws.Range("M27").Value = -993
$ws.Range("N27").Value = -4473.6

$ws.Range("H40").Value = 5779.5713
$ws.Range("I40").Value = 5779.5713
$ws.Range("K40").Value = 5779.5713
$ws.Range("M40").Value = -5643.5713

$ws.Range("H68").Value = 2454478.5
$ws.Range("I68").Value = 4631038
$ws.Range("K68").Value = 4631038
$ws.Range("M68").Value = -4630289

$ws.Range("H71").Value = 2454478.5
$ws.Range("I71").Value = 4631038
$ws.Range("K71").Value = 23155190
$ws.Range("M71").Value = -23151446

$ws.Range("H122").Value = 4128.9756
$ws.Range("J122").Value = 9998.25
$ws.Range("L122").Value = 29994.75
$ws.Range("N122").Value = -34894.75

$ws.Range("H136").Value = 5807.364
$ws.Range("I136").Value = 1897.3334
$ws.Range("J136").Value = 10499.4
$ws.Range("K136").Value = 5692.0002
$ws.Range("L136").Value = 31498.2
$ws.Range("M136").Value = -3142.0002
$ws.Range("N136").Value = -36598.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 315629.34
$ws.Range("I132").Value = 3416.25
$ws.Range("J132").Value = 1252268.6
$ws.Range("K132").Value = 10248.75
$ws.Range("L132").Value = 3756805.8
$ws.Range("M132").Value = -7718.75
$ws.Range("N132").Value = -3761865.8

$ws.Range("H136").Value = 224021.05
$ws.Range("I136").Value = 1697.7179
$ws.Range("J136").Value = 1669122.6
$ws.Range("K136").Value = 5093.153700000001
$ws.Range("L136").Value = 5007367.800000001
$ws.Range("M136").Value = -2543.153700000001
$ws.Range("N136").Value = -5007367.800000001
